# Update "paises.xlsx" COVID dashboard: refresh country statistics and the
# "last updated" timestamp. (commit: "Update countries & provincias Spain")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados ..." timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Julio de 2020 a las 23:23"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 3886132
$ws.Range("C4").Value = 52861
$ws.Range("D4").Value = 1797738
$ws.Range("E4").Value = 1945169
$ws.Range("G4").Value = 348
$ws.Range("H4").Value = 143225

# --- Egipto (row 27) ---
$ws.Range("B27").Value = 87775
$ws.Range("C27").Value = 603
$ws.Range("D27").Value = 28380
$ws.Range("E27").Value = 55093
$ws.Range("G27").Value = 51
$ws.Range("H27").Value = 4302

# --- Israel (row 44) ---
$ws.Range("B44").Value = 50289
$ws.Range("C44").Value = 924
$ws.Range("D44").Value = 21675
$ws.Range("E44").Value = 28205

# --- Barein overtakes Nigeria in the ranking (rows 50 & 51 swap) ---
# Row 50 becomes Barein (updated figures), row 51 becomes Nigeria (unchanged figures).
$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 36422
$ws.Range("C50").Value = 418
$ws.Range("D50").Value = 32372
$ws.Range("E50").Value = 3924
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 126

$ws.Range("A51").Value = "Nigeria"
$ws.Range("B51").Value = 36107
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 14938
$ws.Range("E51").Value = 20391
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 778

# --- Guinea (row 93) ---
$ws.Range("B93").Value = 6544
$ws.Range("C93").Value = 53
$ws.Range("D93").Value = 5511
$ws.Range("E93").Value = 994

# --- Barbados (row 185) ---
$ws.Range("B185").Value = 105
$ws.Range("C185").Value = 1
$ws.Range("E185").Value = 7
